$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows (43-45) of function reference data, matching
# the existing table layout: column A = function name, B = section
# number, C = purpose/note.

$ws.Range("A43").Value = "map()"
$ws.Range("A44").Value = "map_lgl"
$ws.Range("C43").Value = "Returns a list the same length as .x."
$ws.Range("C44").Value = " returns a logical vector"
$ws.Range("A45").Value = "which()"
$ws.Range("B45").Value = "4.8.3"

# Match the saved view state from the target workbook.
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("C32").Select()
